# Generate Report for Handback
# Updates the handback-status report with a new run's data:
#  - file #1's UUID  d6daab99-a7b2-49d3-934c-621dbeee224a -> 8f43d953-10ec-4737-b24c-d3f73de2a9de
#  - file #2's UUID  f4c4224a-7307-4a7e-88f5-81b7a6155c24 -> ffff230b94f4-eebc-42d4-953d-7056c4970e00
#  - refreshed xliff hash / timestamps
# Note: the hyperlink *targets* (Address) keep pointing at the original
# commit's file names (the .rels relationships are untouched by this
# change) -- only the displayed link text is refreshed to the new name.

$wb = $excel.ActiveWorkbook

$uuid1New = "8f43d953-10ec-4737-b24c-d3f73de2a9de"
$uuid2New = "ffff230b94f4-eebc-42d4-953d-7056c4970e00"
$xlfZhCn  = "8f43d953-10ec-4737-b24c-d3f73de2a9de.2c0a9d68f5e89e34d7d5e4983a16db30d1c6744e.zh-cn.xlf"
$xlfDeDe  = "8f43d953-10ec-4737-b24c-d3f73de2a9de.2c0a9d68f5e89e34d7d5e4983a16db30d1c6744e.de-de.xlf"

$dateOverview   = "2016-08-13 09:16:08"
$dateZhHandoff  = "2016-08-13 09:15:57"
$dateZhHandback = "2016-08-13 09:16:26"
$dateDeHandback = "2016-08-13 09:16:36"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$uuid1New.md"
$ws.Range("G2").Value = $dateOverview

$r = $ws.Range("B2")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/d6daab99-a7b2-49d3-934c-621dbeee224a.md", [Type]::Missing, [Type]::Missing, "e2e\$uuid1New.md")

$ws.Range("A3").Value = "$uuid2New.md"
$ws.Range("G3").Value = $dateOverview

$r = $ws.Range("B3")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/f4c4224a-7307-4a7e-88f5-81b7a6155c24.md", [Type]::Missing, [Type]::Missing, "e2e\$uuid2New.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$uuid1New.md"
$ws.Range("I2").Value = "$uuid1New.md"
$ws.Range("G2").Value = $xlfZhCn
$ws.Range("H2").Value = $dateZhHandoff
$ws.Range("J2").Value = $xlfZhCn
$ws.Range("K2").Value = $dateZhHandback

$r = $ws.Range("A2")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/d6daab99-a7b2-49d3-934c-621dbeee224a.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md")

$r = $ws.Range("I2")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e558644688061292ee520e3e0861f3879b42ba3a/e2e/d6daab99-a7b2-49d3-934c-621dbeee224a.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md")

$ws.Range("A3").Value = "$uuid2New.md"
$ws.Range("I3").Value = "$uuid2New.md"
$ws.Range("G3").Value = $xlfZhCn
$ws.Range("H3").Value = $dateZhHandoff
$ws.Range("J3").Value = $xlfZhCn
$ws.Range("K3").Value = $dateZhHandback

$r = $ws.Range("A3")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/f4c4224a-7307-4a7e-88f5-81b7a6155c24.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md")

$r = $ws.Range("I3")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e558644688061292ee520e3e0861f3879b42ba3a/e2e/f4c4224a-7307-4a7e-88f5-81b7a6155c24.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$uuid1New.md"
$ws.Range("I2").Value = "$uuid1New.md"
$ws.Range("G2").Value = $xlfDeDe
$ws.Range("H2").Value = $dateOverview
$ws.Range("J2").Value = $xlfDeDe
$ws.Range("K2").Value = $dateDeHandback

$r = $ws.Range("A2")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/d6daab99-a7b2-49d3-934c-621dbeee224a.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md")

$r = $ws.Range("I2")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b51564eca177d3838be4b36505b88eedc8f43d80/e2e/d6daab99-a7b2-49d3-934c-621dbeee224a.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md")

$ws.Range("A3").Value = "$uuid2New.md"
$ws.Range("I3").Value = "$uuid2New.md"
$ws.Range("G3").Value = $xlfDeDe
$ws.Range("H3").Value = $dateOverview
$ws.Range("J3").Value = $xlfDeDe
$ws.Range("K3").Value = $dateDeHandback

$r = $ws.Range("A3")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/oltest/blob/cbf08d233e120eca98a5ba7299848d96b34d1e77/e2e/f4c4224a-7307-4a7e-88f5-81b7a6155c24.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md")

$r = $ws.Range("I3")
$r.Hyperlinks.Delete()
$ws.Hyperlinks.Add($r, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b51564eca177d3838be4b36505b88eedc8f43d80/e2e/f4c4224a-7307-4a7e-88f5-81b7a6155c24.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md")
